# edit.ps1 -- apply LOM3098 syllabus content restructuring
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix Objetivos body text (row 10, B/C already carry the right style) ---
$ws.Range("B10").Value2 = 'Propiciar conhecimentos teóricos e práticos de materiais e ferramentas abrasivas, incluindo as matérias primas, os processos de fabricação, qualificação, teste e aplicação de ferramentas abrasivas; e os mecanismos de abrasão atuantes durante os processos de usinagem. Aprimorar a formação do Engenheiro de Materiais numa área estratégica de Ciência e Engenharia de Materiais presente em diversos segmentos da indústria de transformação do país'
$ws.Range("C10").Value2 = 'Propiciar conhecimentos teóricos e práticos de materiais e ferramentas abrasivas, incluindo as matérias primas, os processos de fabricação, qualificação, teste e aplicação de ferramentas abrasivas; e os mecanismos de abrasão atuantes durante os processos de usinagem. Aprimorar a formação do Engenheiro de Materiais numa área estratégica de Ciência e Engenharia de Materiais presente em diversos segmentos da indústria de transformação do país'

# --- Step 2: insert two blank rows at 13:14 to make room for the extra "Docentes responsaveis" entry ---
$ws.Range("A13:A14").EntireRow.Insert()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# --- Step 3: rewrite rows 13-25 content, styles and heights to match target layout ---
# Row 13
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value2 = '519033 - Carlos Yujiro Shigue'
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value2 = '519033 - Carlos Yujiro Shigue'
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value2 = '5817692 - Katia Cristiane Gandolpho Candioto'
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value2 = '5817692 - Katia Cristiane Gandolpho Candioto'
$ws.Rows.Item(14).AutoFit()

# Row 15
$ws.Range("A3").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value2 = 'Programa resumido:'
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value2 = 'Introdução aos materiais e ferramentas abrasivas. Matérias primas utilizadas na fabricação de ferramentas abrasivas. Processos de fabricação de ferramentas abrasivas. Caracterização, teste e inspeção de ferramentas abrasivas. Mecânica da usinagem com ferramentas abrasivas. Avaliação de desempenho. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.'
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value2 = 'Introdução aos materiais e ferramentas abrasivas. Matérias primas utilizadas na fabricação de ferramentas abrasivas. Processos de fabricação de ferramentas abrasivas. Caracterização, teste e inspeção de ferramentas abrasivas. Mecânica da usinagem com ferramentas abrasivas. Avaliação de desempenho. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A3").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value2 = 'Short syllabus:'
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value2 = 'Programa:'
$ws.Range("B10").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value2 = 'Conteúdo teórico: 1. Introdução aos materiais e ferramentas abrasivas: histórico, materiais abrasivos naturais e sintéticos e características principais dos materiais abrasivos. Dados econômicos das ferramentas abrasivas.2. Matérias primas utilizadas na fabricação de ferramentas abrasivas: cerâmicas, borrachas, metais e polímeros. Processos de obtenção das matérias primas.3. Processos de fabricação de ferramentas abrasivas: discos, rebolos, pontas montadas e lixas.4. Caracterização, teste e inspeção de ferramentas abrasivas: ensaios destrutivos e não destrutivos. Normas e códigos de segurança. 5. Mecânica da usinagem com ferramentas abrasivas. Operações com abrasivos: corte, retificação, desbaste, acabamento, lapidação e afiação.6. Avaliação de desempenho: aspectos térmicos, refrigeração, lubrificação, rugosidade superficial, interação metal-ferramenta e defeitos em ferramentas abrasivas.Conteúdo prático: 1. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.2. Visita a fabricantes de abrasivos. 3. Visita a usuários de ferramentas abrasivas.'
$ws.Range("C10").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value2 = 'Conteúdo teórico: 1. Introdução aos materiais e ferramentas abrasivas: histórico, materiais abrasivos naturais e sintéticos e características principais dos materiais abrasivos. Dados econômicos das ferramentas abrasivas.2. Matérias primas utilizadas na fabricação de ferramentas abrasivas: cerâmicas, borrachas, metais e polímeros. Processos de obtenção das matérias primas.3. Processos de fabricação de ferramentas abrasivas: discos, rebolos, pontas montadas e lixas.4. Caracterização, teste e inspeção de ferramentas abrasivas: ensaios destrutivos e não destrutivos. Normas e códigos de segurança. 5. Mecânica da usinagem com ferramentas abrasivas. Operações com abrasivos: corte, retificação, desbaste, acabamento, lapidação e afiação.6. Avaliação de desempenho: aspectos térmicos, refrigeração, lubrificação, rugosidade superficial, interação metal-ferramenta e defeitos em ferramentas abrasivas.Conteúdo prático: 1. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.2. Visita a fabricantes de abrasivos. 3. Visita a usuários de ferramentas abrasivas.'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value2 = 'Syllabus:'
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Range("A3").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value2 = 'Avaliação:'
$ws.Rows.Item(19).AutoFit()

# Row 20
$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value2 = 'Método:'
$ws.Range("B10").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value2 = 'A avaliação será constituída por aulas expositivas, aulas de exercícios e práticas laboratoriais. Serão aplicadas pelo menos duas avaliações.'
$ws.Range("C10").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value2 = 'A avaliação será constituída por aulas expositivas, aulas de exercícios e práticas laboratoriais. Serão aplicadas pelo menos duas avaliações.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value2 = 'Critério:'
$ws.Range("B10").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value2 = 'A nota final será a média das avaliações escritas e práticas'
$ws.Range("C10").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value2 = 'A nota final será a média das avaliações escritas e práticas'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A3").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value2 = 'Norma de recuperação:'
$ws.Range("B10").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value2 = 'A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2.'
$ws.Range("C10").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value2 = 'A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2.'
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value2 = 'Bibliografia:'
$ws.Range("B10").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value2 = '1. NUSSBAUM, G. C. Rebolos e abrasivos. Tecnologia básica. São Paulo: Ícone Editora, 1988. 2. KLOCKE, F. Manufacturing processes 2. Grinding, honing, lapping. Berlim: Springer Verlag, 2009.3. MALKIN, S.; GUO, C. Grinding technology: theory and application of machining with abrasives. New York: Industrial Press Inc., 2008.4. JACKSON, M. J.; DAVIM, J. P. Machining with abrasives. New York: Springer Science, 2011.5. FERRARESI, D. Usinagem dos metais. São Paulo: Editora Edgard Blucher, 1970.6. STEMMER, C. E. Ferramentas de corte II: brocas, alargadores, ferramentas de rocar, fresas, brochas, rebolos e abrasivos. Florianópolis: Editora da UFSC, 1995.7. KINGERY, W. D. Ceramic fabrication process. New York: John Wiley, 1958.8. GARDZIELLA, A.; PILATO, L.A.; KNOP, A. Phenolic resins: chemistry, applications, standardization, safety and ecology. Berlim: Springer Verlag, 2000.9. MARINESCU, Ioan D. Tribology of abrasive machining processes. 2ª Ed. New York: Willian Andrew, 2004.'
$ws.Range("C10").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value2 = '1. NUSSBAUM, G. C. Rebolos e abrasivos. Tecnologia básica. São Paulo: Ícone Editora, 1988. 2. KLOCKE, F. Manufacturing processes 2. Grinding, honing, lapping. Berlim: Springer Verlag, 2009.3. MALKIN, S.; GUO, C. Grinding technology: theory and application of machining with abrasives. New York: Industrial Press Inc., 2008.4. JACKSON, M. J.; DAVIM, J. P. Machining with abrasives. New York: Springer Science, 2011.5. FERRARESI, D. Usinagem dos metais. São Paulo: Editora Edgard Blucher, 1970.6. STEMMER, C. E. Ferramentas de corte II: brocas, alargadores, ferramentas de rocar, fresas, brochas, rebolos e abrasivos. Florianópolis: Editora da UFSC, 1995.7. KINGERY, W. D. Ceramic fabrication process. New York: John Wiley, 1958.8. GARDZIELLA, A.; PILATO, L.A.; KNOP, A. Phenolic resins: chemistry, applications, standardization, safety and ecology. Berlim: Springer Verlag, 2000.9. MARINESCU, Ioan D. Tribology of abrasive machining processes. 2ª Ed. New York: Willian Andrew, 2004.'
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$ws.Range("A3").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value2 = 'Requisitos:'
$ws.Rows.Item(24).AutoFit()

# Row 25
$ws.Range("B10").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value2 = 'LOM3011 -  Ensaios Mecânicos  (Requisito)
'
$ws.Range("C10").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value2 = 'LOM3011 -  Ensaios Mecânicos  (Requisito)
'
$ws.Rows.Item(25).RowHeight = 30

$excel.CutCopyMode = $false